# CoderDojoKells Autumn Session 8 - 20181121 update
#
# This edit reproduces a LibreOffice-style EMU/1-100mm round-trip on every
# shape's size (width & height each shrink by 360 EMU == 1/100 mm), plus one
# shape (the flipped background accent on slide master 1) also shifting left
# by 720 EMU. It also removes the "Parents, pitch in" bullet from slide 4 and
# hides slide 8.

$p = $ppt.ActivePresentation

# --- helpers -----------------------------------------------------------
# Converting EMU -> points loses precision (points are stored/truncated as
# float32 internally), so add a small sub-EMU epsilon before dividing to make
# sure the value survives the round trip and lands back on the exact EMU we
# want.
function EmuToPt($emu) {
    return ($emu + 0.45) / 12700.0
}

function ShrinkShapeExtent($shp, $deltaEmu) {
    $curW = [Math]::Round($shp.Width * 12700)
    $curH = [Math]::Round($shp.Height * 12700)
    $shp.Width  = EmuToPt($curW - $deltaEmu)
    $shp.Height = EmuToPt($curH - $deltaEmu)
}

# --- 1. every shape on every slide loses 360 EMU off width & height ----
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)
        ShrinkShapeExtent $shp 360
    }
}

# --- 2. every non-placeholder shape on every slide master loses 360 EMU
#        off width & height too; the flipped custom-geometry shape on the
#        first master also shifts 720 EMU to the left -------------------
for ($d = 1; $d -le $p.Designs.Count; $d++) {
    $design = $p.Designs.Item($d)
    $master = $design.SlideMaster
    for ($j = 1; $j -le $master.Shapes.Count; $j++) {
        $shp = $master.Shapes.Item($j)
        if ($shp.Type -eq 14) {
            # msoPlaceholder - title/body placeholders are untouched
            continue
        }
        ShrinkShapeExtent $shp 360
        if ($shp.HorizontalFlip -and $shp.VerticalFlip) {
            $curLeft = [Math]::Round($shp.Left * 12700)
            $shp.Left = EmuToPt($curLeft - 720)
        }
    }
}

# --- 3. slide 4: drop the "Parents, pitch in" bullet --------------------
$s4 = $p.Slides.Item(4)
$motto = $s4.Shapes.Item(2)
$tr = $motto.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count
for ($k = $paraCount; $k -ge 1; $k--) {
    $para = $tr.Paragraphs($k, 1)
    if ($para.Text -match "Parents, pitch in") {
        $para.Delete()
    }
}

# --- 4. slide 8 is hidden from the slide show ---------------------------
$s8 = $p.Slides.Item(8)
$s8.SlideShowTransition.Hidden = $true
